$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 888-889), which
# shifts the existing rows 888-1011 down to 890-1013 while preserving all
# values and formatting.
$ws.Rows.Item(888).Resize(2).Insert()

# Populate the new row 888 ("Primera") with the new weekly record.
$ws.Range("A888").Value = 8
$ws.Range("B888").Value = "Terminal La Palmera de La Serena"
$ws.Range("C888").Value = "Coquimbo"
$ws.Range("D888").Value = 44984
$ws.Range("E888").Value = 4
$ws.Range("F888").Value = 100112008
$ws.Range("G888").Value = "Coliflor"
$ws.Range("H888").Value = "Sin especificar"
$ws.Range("I888").Value = "Primera"
$ws.Range("J888").Value = 2000
$ws.Range("K888").Value = 700
$ws.Range("L888").Value = 800
$ws.Range("M888").Value = 750
$ws.Range("N888").Value = "`$/unidad"
$ws.Range("O888").Value = "Provincia del Elquí"
$ws.Range("P888").Value = 750
$ws.Range("Q888").Value = 1
$ws.Range("R888").Value = "Hortaliza"

# Populate the new row 889 ("Segunda") with the new weekly record.
$ws.Range("A889").Value = 8
$ws.Range("B889").Value = "Terminal La Palmera de La Serena"
$ws.Range("C889").Value = "Coquimbo"
$ws.Range("D889").Value = 44984
$ws.Range("E889").Value = 4
$ws.Range("F889").Value = 100112008
$ws.Range("G889").Value = "Coliflor"
$ws.Range("H889").Value = "Sin especificar"
$ws.Range("I889").Value = "Segunda"
$ws.Range("J889").Value = 1360
$ws.Range("K889").Value = 500
$ws.Range("L889").Value = 600
$ws.Range("M889").Value = 550
$ws.Range("N889").Value = "`$/unidad"
$ws.Range("O889").Value = "Provincia del Elquí"
$ws.Range("P889").Value = 550
$ws.Range("Q889").Value = 1
$ws.Range("R889").Value = "Hortaliza"

# Make sure the date cells keep the expected date number format.
$ws.Range("D888").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D889").NumberFormat = "YYYY-MM-DD HH:MM:SS"
